$d = $word.ActiveDocument

# 1. Merge "pro" + "jeto" -> "projeto" (simple text fix, no visible change but normalizes runs)
$d.Content.Find.Execute("necessários para o projeto", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "necessários para o projeto", 2)

# 2. "implementou" -> "desenvolveu"
$d.Content.Find.Execute("implementou um novo c", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "desenvolveu um novo c", 2)

# 3. "a o Arduino" -> "para o Arduino"
$d.Content.Find.Execute("digo a o Arduino", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "digo para o Arduino", 2)

# 4. "Wi-fi.  " -> "Wi-fi co sucesso.\nOBS: ..."
$d.Content.Find.Execute("rede Wi-fi.  ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "rede Wi-fi co sucesso.^lOBS: Não foi possível se conectar a redes que possuem acentos em seus SSIDs, também não foi possível se conectar a redes de 5GHz  ", 2)
